$d = $word.ActiveDocument

# --- Part 1: relocate the "_GoBack" bookmark in paragraph 1 so that it
#     ends up after the "-F-AutoIndentationTest-U-D-W" and "测试用例文档"
#     runs instead of before them (visible text/formatting is unchanged;
#     only the bookmark position moves to the end of the paragraph's text).

$p1 = $d.Paragraphs(1).Range
$p1TextEnd = $p1.End - 1   # position right after the last visible character, before the paragraph mark

$d.Bookmarks("_GoBack").Delete()

# Work around an edge-case in zero-length bookmark placement at the very
# last text position of a paragraph: temporarily insert a marker character
# after that position, anchor the bookmark there (no longer the "last"
# position), then remove the marker again. The bookmark collapses back to
# a zero-length bookmark in the correct spot.
$marker = $d.Range($p1TextEnd, $p1TextEnd)
$marker.InsertAfter("X")

$bmRange = $d.Range($p1TextEnd, $p1TextEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($p1TextEnd, $p1TextEnd + 1).Delete()

# --- Part 2: expand the test-target sentence in paragraph 2.
$d.Content.Find.Execute("当大括号中有内容时，能否实现自动缩进", $true, $false, $false, $false, $false,
                         $true, 1, $false, "当大括号中有内容，且与左括号位于同一行时，能否实现自动缩进", 2)
